$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---------------------------------------------------------------
# Before: B2:C2 merged, E2:G2 merged (C2 and F2 are non-anchor duplicates)
# After:  B2:C2 stays merged (C2 cleared); E2:G2 splits into three 1-cell
#         merges E2:E2, F2:F2, G2:G2 (F2 cleared)
$ws.Range("B2:C2").UnMerge()
$ws.Range("C2").Value = $null
$ws.Range("B2:C2").Merge()

$ws.Range("E2:G2").UnMerge()
$ws.Range("F2").Value = $null
$ws.Range("E2:E2").Merge()
$ws.Range("F2:F2").Merge()
$ws.Range("G2:G2").Merge()

# --- Row 3 ---------------------------------------------------------------
# Before: C3:D3 merged, E3:G3 merged, H3:I3 merged
# After:  same merges, but D3 / E3 / F3 / G3 / I3 cleared (already-blank or
#         non-anchor cells disappear from the sheet data)
$ws.Range("C3:D3").UnMerge()
$ws.Range("D3").Value = $null
$ws.Range("C3:D3").Merge()

$ws.Range("E3:G3").UnMerge()
$ws.Range("E3").Value = $null
$ws.Range("F3").Value = $null
$ws.Range("G3").Value = $null
$ws.Range("E3:G3").Merge()

$ws.Range("H3:I3").UnMerge()
$ws.Range("I3").Value = $null
$ws.Range("H3:I3").Merge()

# --- Row 4 ---------------------------------------------------------------
# Before: B4:C4 merged, E4:F4 merged, H4:I4 merged
# After:  same merges, but C4 / E4 / F4 / I4 cleared
$ws.Range("B4:C4").UnMerge()
$ws.Range("C4").Value = $null
$ws.Range("B4:C4").Merge()

$ws.Range("E4:F4").UnMerge()
$ws.Range("E4").Value = $null
$ws.Range("F4").Value = $null
$ws.Range("E4:F4").Merge()

$ws.Range("H4:I4").UnMerge()
$ws.Range("I4").Value = $null
$ws.Range("H4:I4").Merge()

# --- Row 5 ---------------------------------------------------------------
# Before: B5:D5 merged, H5:I5 merged
# After:  B5:D5 splits into B5:C5 + D5:D5 (C5 / D5 cleared); H5:I5 stays
#         merged but I5 cleared
$ws.Range("B5:D5").UnMerge()
$ws.Range("C5").Value = $null
$ws.Range("D5").Value = $null
$ws.Range("B5:C5").Merge()
$ws.Range("D5:D5").Merge()

$ws.Range("H5:I5").UnMerge()
$ws.Range("I5").Value = $null
$ws.Range("H5:I5").Merge()

# --- Row 6 ---------------------------------------------------------------
# Before: C6:D6 merged, H6:I6 merged (both stay the same shape)
# After:  D6 and I6 (non-anchor duplicates) are blanked out
$ws.Range("C6:D6").UnMerge()
$ws.Range("D6").Value = $null
$ws.Range("C6:D6").Merge()

$ws.Range("H6:I6").UnMerge()
$ws.Range("I6").Value = $null
$ws.Range("H6:I6").Merge()
